# Update NASDAQ-100 ticker list (data as of 2024-03-24):
#  - Insert "Linde plc" (LIN) as a new constituent, placed alphabetically
#    right before "Lululemon" (LULU), which pushes all following rows down.
#  - Remove "Splunk" (SPLK), which is no longer a standalone constituent
#    (acquired by Cisco), restoring the original row count.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 59 (current LULU row); everything from
# row 59 down (including the SPLK row at 88) shifts down by one.
$ws.Rows("59:59").Insert()

# Populate the newly inserted row 59 with the Linde plc data.
$ws.Range("A59").Value = "LIN"
$ws.Range("B59").Value = "Linde plc"
$ws.Range("C59").Value = "Materials"
$ws.Range("D59").Value = "Industrial Gases"

# SPLK (originally row 88) is now at row 89 after the insert above.
# Delete that row entirely so the table returns to its original length.
$ws.Rows("89:89").Delete()
